# T_SOL_summary_NonMast: refresh the seasonal latency/relative-timing table
# with the updated "nest measures" figures (rows stay Description/Latency to
# AB/AB Rel. to Sunrise/Latency to QB/QB Rel. to Sunset; columns stay
# Winter/Spring/Summer/Autumn/All). Only the B2:F5 data block changes value;
# headers (row 1) and row labels (column A) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Latency to AB
$ws.Range("B2").Value = "44.14 ± 27.69"
$ws.Range("C2").Value = "42.88 ± 28.84"
$ws.Range("D2").Value = "26.99 ± 15.00"
$ws.Range("E2").Value = "27.67 ± 15.44"
$ws.Range("F2").Value = "40.91 ± 26.87"

# Row 3 - AB Rel. to Sunrise
$ws.Range("B3").Value = "13.93 ± 70.73"
$ws.Range("C3").Value = "44.43 ± 78.70"
$ws.Range("D3").Value = "37.23 ± 50.16"
$ws.Range("E3").Value = "-5.30 ± 23.45"
$ws.Range("F3").Value = "17.45 ± 68.63"

# Row 4 - Latency to QB
$ws.Range("B4").Value = "74.90 ± 101.84"
$ws.Range("C4").Value = "107.13 ± 116.94"
$ws.Range("D4").Value = "108.07 ± 129.21"
$ws.Range("E4").Value = "64.78 ± 78.02"
$ws.Range("F4").Value = "80.42 ± 103.67"

# Row 5 - QB Rel. to Sunset
$ws.Range("B5").Value = "-66.68 ± 90.98"
$ws.Range("C5").Value = "-108.41 ± 101.03"
$ws.Range("D5").Value = "-63.19 ± 72.09"
$ws.Range("E5").Value = "-31.65 ± 56.45"
$ws.Range("F5").Value = "-69.07 ± 90.94"

# Columns B, D and F got a touch wider in the refreshed export.
$ws.Columns.Item(2).ColumnWidth = 12.25
$ws.Columns.Item(4).ColumnWidth = 13.25
$ws.Columns.Item(6).ColumnWidth = 12.25
